# CryCompanywiseStockReport_1.xlsx - stock quantity/value correction pass.
#
# The sheet is a flat "Companywise Stock" report with columns:
#   A=SlNo  B=Code  C=Item  D=Rate  E=OtherRate  F=Qty  G=Value(=D*Qty)
# Each company block ends in a "Sub Total:" row whose B cell is the sum of
# the G values in that block; row 619/620 hold the overall Sub Total/Grand
# Total, i.e. the sum of every block's Sub Total.
#
# This edit corrects the Qty (F) for a number of line items (mostly -1/-2
# adjustments, a couple of larger corrections), recomputes the matching
# Value (G = Rate * Qty), fixes two pairs of rows (227/228 and 243/244)
# where the Code/OtherRate/Qty/Value had been swapped between two
# identically-named items, and finally rolls all of that up through the
# affected company Sub Totals and the final Sub Total / Grand Total cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# --- Line-item Qty/Value corrections (G = D * new Qty) ------------------
Set-Cell "F6"   52
Set-Cell "G6"   1553.76

Set-Cell "F61"  48
Set-Cell "G61"  3371.52

Set-Cell "F71"  298
Set-Cell "G71"  18982.6

Set-Cell "F77"  237
Set-Cell "G77"  11077.38

Set-Cell "F78"  37
Set-Cell "G78"  2105.3

Set-Cell "F84"  18
Set-Cell "G84"  1844.28

Set-Cell "F115" 177
Set-Cell "G115" 17135.37

Set-Cell "F135" 20
Set-Cell "G135" 620.6

Set-Cell "F141" 43
Set-Cell "G141" 2301.79

Set-Cell "F146" 14
Set-Cell "G146" 1178.66

Set-Cell "F149" 212
Set-Cell "G149" 13737.6

Set-Cell "F150" 23
Set-Cell "G150" 1069.27

Set-Cell "F167" 12
Set-Cell "G167" 3444.12

Set-Cell "F198" 20
Set-Cell "G198" 5442.4

Set-Cell "F203" 50
Set-Cell "G203" 1008

Set-Cell "F225" 71
Set-Cell "G225" 8110.33

# Rows 227/228: Code (B), OtherRate (E), Qty (F) and Value (G) were
# swapped between these two "HUL-Kissan nango jam 490g" rows.
Set-Cell "B227" 55373
Set-Cell "E227" 163.62
Set-Cell "F227" -94
Set-Cell "G227" -13562.32

Set-Cell "B228" 63520
Set-Cell "E228" 153.4
Set-Cell "F228" 65
Set-Cell "G228" 9378.200000000001

# Rows 243/244: same swap, for the "Hul-pears pure and gentle 3x125 gm" rows.
Set-Cell "B243" 63560
Set-Cell "E243" 134.87
Set-Cell "F243" 1
Set-Cell "G243" 126.86

Set-Cell "B244" 60325
Set-Cell "E244" 151.57
Set-Cell "F244" -102
Set-Cell "G244" -12939.72

Set-Cell "F247" 129
Set-Cell "G247" 13404.39

Set-Cell "F280" 128
Set-Cell "G280" 21649.92

Set-Cell "F291" 102
Set-Cell "G291" 4387.02

Set-Cell "F294" 24
Set-Cell "G294" 1712.64

Set-Cell "F296" 31
Set-Cell "G296" 657.2

Set-Cell "F300" 161
Set-Cell "G300" 20110.51

Set-Cell "F302" 28
Set-Cell "G302" 5904.92

Set-Cell "F303" 20
Set-Cell "G303" 4217.8

Set-Cell "F320" 36
Set-Cell "G320" 2471.4

Set-Cell "F326" 58
Set-Cell "G326" 1724.92

Set-Cell "F333" 39
Set-Cell "G333" 1945.71

Set-Cell "F334" 187
Set-Cell "G334" 9690.34

Set-Cell "F338" 72
Set-Cell "G338" 1706.4

Set-Cell "F343" 29
Set-Cell "G343" 2087.13

Set-Cell "F345" 30
Set-Cell "G345" 1842.3

Set-Cell "F454" 45
Set-Cell "G454" 1536.75

Set-Cell "F498" 0
Set-Cell "G498" 0

Set-Cell "F499" 0
Set-Cell "G499" 0

Set-Cell "F509" 190
Set-Cell "G509" 15272.2

Set-Cell "F549" 22
Set-Cell "G549" 1052.92

Set-Cell "F552" 9
Set-Cell "G552" 916.11

Set-Cell "F555" 13
Set-Cell "G555" 904.28

Set-Cell "F556" 0
Set-Cell "G556" 0

Set-Cell "F577" 37
Set-Cell "G577" 1590.63

Set-Cell "F578" 52
Set-Cell "G578" 2594.28

Set-Cell "F582" 21
Set-Cell "G582" 1196.79

Set-Cell "F599" 1294
Set-Cell "G599" 211064.34

Set-Cell "F601" 357
Set-Cell "G601" 100984.59

Set-Cell "F602" 309
Set-Cell "G602" 44696.85

Set-Cell "F612" 28
Set-Cell "G612" 1147.72

Set-Cell "F613" 125
Set-Cell "G613" 19895

# --- Company "Sub Total:" rows recomputed from the new item values ------
Set-Cell "B10"  26344.56
Set-Cell "B90"  166118.01
Set-Cell "B117" 10976.43
Set-Cell "B138" 2226.68
Set-Cell "B142" 2813.47
Set-Cell "B147" 12290.3
Set-Cell "B156" 28711.59
Set-Cell "B175" 26450.71
Set-Cell "B216" 32083.57
Set-Cell "B260" 166203.76
Set-Cell "B304" 161449
Set-Cell "B330" 25158.29
Set-Cell "B346" 23338.3
Set-Cell "B460" 12074.97
Set-Cell "B500" 0
Set-Cell "B510" 20677.08
Set-Cell "B560" 3143.33
Set-Cell "B583" 12254.77
Set-Cell "B606" 357593.83
Set-Cell "B618" 41112.52

# --- Overall Sub Total / Grand Total -------------------------------------
Set-Cell "B619" 1556534.02
Set-Cell "B620" 1556534.02
